$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1-2 headers (unchanged values, rewritten defensively)
$ws.Range("B1").Value = "s"
$ws.Range("C1").Value = "s"
$ws.Range("A2").Value = "array"
$ws.Range("B2").Value = "index"
$ws.Range("C2").Value = "soundPath"

# Data rows 3-14: narration / tutorial entries
$ws.Range("B3").Value  = "intro_01"
$ws.Range("C3").Value  = "intro_01"
$ws.Range("B4").Value  = "intro_02"
$ws.Range("C4").Value  = "intro_02"
$ws.Range("B5").Value  = "tutorial01"
$ws.Range("C5").Value  = "tutorial01"
$ws.Range("B6").Value  = "tutorial02"
$ws.Range("C6").Value  = "tutorial02"
$ws.Range("B7").Value  = "tutorial03"
$ws.Range("C7").Value  = "tutorial03"
$ws.Range("B8").Value  = "tutorial04"
$ws.Range("C8").Value  = "tutorial04"
$ws.Range("B9").Value  = "tutorial05"
$ws.Range("C9").Value  = "tutorial05"
$ws.Range("B10").Value = "tutorial06"
$ws.Range("C10").Value = "tutorial06"
$ws.Range("B11").Value = "tutorial07"
$ws.Range("C11").Value = "tutorial07"
$ws.Range("B12").Value = "tutorial08"
$ws.Range("C12").Value = "tutorial08"
$ws.Range("B13").Value = "tutorial09"
$ws.Range("C13").Value = "tutorial09"
$ws.Range("B14").Value = "tutorial10"
$ws.Range("C14").Value = "tutorial10"

# Data rows 15-21: new "foreset" quiz-check entries reusing tutorial sound paths
$ws.Range("B15").Value = "foreset_01"
$ws.Range("C15").Value = "tutorial01"
$ws.Range("B16").Value = "foreset_02"
$ws.Range("C16").Value = "tutorial02"
$ws.Range("B17").Value = "foreset_03"
$ws.Range("C17").Value = "tutorial03"
$ws.Range("B18").Value = "foreset_04"
$ws.Range("C18").Value = "tutorial04"
$ws.Range("B19").Value = "foreset_05"
$ws.Range("C19").Value = "tutorial05"
$ws.Range("B20").Value = "foreset_06"
$ws.Range("C20").Value = "tutorial06"
$ws.Range("B21").Value = "foreset_07"
$ws.Range("C21").Value = "tutorial07"

# Restore the selection that Excel recorded at save time
$ws.Range("G17").Select()
